$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'silicone knee pads'
$ws.Cells.Item(2, 1).Value = 'skins compression men'
$ws.Cells.Item(3, 1).Value = 'sliding knee pads baseball'
$ws.Cells.Item(4, 1).Value = 'small football knee pads'
$ws.Cells.Item(5, 1).Value = 'snowboarding knee pad'
$ws.Cells.Item(6, 1).Value = 'soccer knee pad'
$ws.Cells.Item(7, 1).Value = 'soccer knee protector'
$ws.Cells.Item(8, 1).Value = 'soccer pants for boys'
$ws.Cells.Item(9, 1).Value = 'spandex pants for men'
$ws.Cells.Item(10, 1).Value = 'spats bjj'
$ws.Cells.Item(11, 1).Value = 'sports basketball pants'
$ws.Cells.Item(12, 1).Value = 'sports leggings men'
$ws.Cells.Item(13, 1).Value = 'squat protector'
$ws.Cells.Item(14, 1).Value = 'strip pants men'
$ws.Cells.Item(15, 1).Value = 'tactical pants with knee pads'
$ws.Cells.Item(16, 1).Value = 'thermal baselayer men'
$ws.Cells.Item(17, 1).Value = 'thermal compression pants men'
$ws.Cells.Item(18, 1).Value = 'thermal leggings men'
$ws.Cells.Item(19, 1).Value = 'thermal winter pants'
$ws.Cells.Item(20, 1).Value = 'thin knee pads'
$ws.Cells.Item(21, 1).Value = 'tight pant'
$ws.Cells.Item(22, 1).Value = 'tights black'
$ws.Cells.Item(23, 1).Value = 'tights for men'
$ws.Cells.Item(24, 1).Value = 'under amour fleece leggings'
$ws.Cells.Item(25, 1).Value = 'under armour leggings'
$ws.Cells.Item(26, 1).Value = 'underarmor thermal pants mens'
$ws.Cells.Item(27, 1).Value = 'volleyball clothes men'
$ws.Cells.Item(28, 1).Value = 'volleyball compression knee pads'
$ws.Cells.Item(29, 1).Value = 'volleyball knee pads adult'
$ws.Cells.Item(30, 1).Value = 'volleyball knee pads youth'
$ws.Cells.Item(31, 1).Value = 'warm compression pants mens'
$ws.Cells.Item(32, 1).Value = 'weightlifting equipment'
$ws.Cells.Item(33, 1).Value = 'white basketball knee pads'
$ws.Cells.Item(34, 1).Value = 'white tights mens basketball'
$ws.Cells.Item(35, 1).Value = 'womens basketball knee pads'
$ws.Cells.Item(36, 1).Value = 'workout pads for hands'
$ws.Cells.Item(37, 1).Value = 'workout squat pad'
$ws.Cells.Item(38, 1).Value = 'wrestling knee pads adult'
$ws.Cells.Item(39, 1).Value = 'wrestling tights boys'
$ws.Cells.Item(40, 1).Value = 'xl volleyball knee pads'
$ws.Cells.Item(41, 1).Value = 'youth basketball knee'
$ws.Cells.Item(42, 1).Value = 'youth basketball leggings for boys'
$ws.Cells.Item(43, 1).Value = 'youth basketball pants girls'
$ws.Cells.Item(44, 1).Value = 'youth football knee pads'
$ws.Cells.Item(45, 1).Value = 'youth knee pads for skating'
$ws.Cells.Item(46, 1).Value = 'youth knee pads mountain bike'
$ws.Cells.Item(47, 1).Value = 'youth soccer gear'
$ws.Cells.Item(48, 1).Value = 'youth sports tights'
$ws.Cells.Item(49, 1).Value = 'mens running tights capri'
$ws.Cells.Item(50, 1).Value = 'mens basketball pants tall'
$ws.Cells.Item(51, 1).Value = 'volleyball hip protectors'
$ws.Cells.Item(52, 1).Value = 'knee compression basketball'
$ws.Cells.Item(53, 1).Value = 'sports compression pants'
$ws.Cells.Item(54, 1).Value = 'basketball leggings youth'
$ws.Cells.Item(55, 1).Value = 'knee pads for running'
$ws.Cells.Item(56, 1).Value = 'wrestling pants for men'
$ws.Cells.Item(57, 1).Value = 'male compression pants'
$ws.Cells.Item(58, 1).Value = 'cheap compression pants men'
$ws.Cells.Item(59, 1).Value = 'men compression pants pack'
$ws.Cells.Item(60, 1).Value = 'knee pad men'
$ws.Cells.Item(61, 1).Value = 'volleyball knee pads for youth'
$ws.Cells.Item(62, 1).Value = 'leggings for man'
$ws.Cells.Item(63, 1).Value = 'volleyball knee pads xl'
$ws.Cells.Item(64, 1).Value = 'spandex leggings men'
$ws.Cells.Item(65, 1).Value = 'knees compression'
$ws.Cells.Item(66, 1).Value = 'athletic boys pants'
$ws.Cells.Item(67, 1).Value = 'soccer knee pads'
$ws.Cells.Item(68, 1).Value = 'adults knee pads'
$ws.Cells.Item(69, 1).Value = 'leg compression leggings'
$ws.Cells.Item(70, 1).Value = 'big and tall tights for men'
$ws.Cells.Item(71, 1).Value = 'wrestling mens apparel'
$ws.Cells.Item(72, 1).Value = 'mens hiking leggings'
$ws.Cells.Item(73, 1).Value = 'softball pants mens'
$ws.Cells.Item(74, 1).Value = 'running compression knee'
$ws.Cells.Item(75, 1).Value = 'compression pads for basketball'
$ws.Cells.Item(76, 1).Value = 'gym tights'
$ws.Cells.Item(77, 1).Value = 'boys soccer leggings'
$ws.Cells.Item(78, 1).Value = 'hockey pants men'
$ws.Cells.Item(79, 1).Value = 'knee pads for volleyball for men'
$ws.Cells.Item(80, 1).Value = 'baseball pants men'
$ws.Cells.Item(81, 1).Value = 'mens running pants'
$ws.Cells.Item(82, 1).Value = 'sports pants for men'
$ws.Cells.Item(83, 1).Value = 'boys snowboarding pants'
$ws.Cells.Item(84, 1).Value = 'clothing protector'
$ws.Cells.Item(85, 1).Value = 'baseball pants for youth'
$ws.Cells.Item(86, 1).Value = 'mens gym pants'
$ws.Cells.Item(87, 1).Value = 'football pants adult'
$ws.Cells.Item(88, 1).Value = 'knee pads sports'
$ws.Cells.Item(89, 1).Value = 'protective knee pads'
$ws.Cells.Item(90, 1).Value = 'girls volleyball pads'
$ws.Cells.Item(91, 1).Value = 'pants youth'
$ws.Cells.Item(92, 1).Value = 'big boys compression leggings'
$ws.Cells.Item(93, 1).Value = 'basketball spandex'
$ws.Cells.Item(94, 1).Value = 'mens basketball clothing'
$ws.Cells.Item(95, 1).Value = 'boys tights for sports'
$ws.Cells.Item(96, 1).Value = 'leggings for cycling'
$ws.Cells.Item(97, 1).Value = 'hiking capri men'
$ws.Cells.Item(98, 1).Value = 'baseball pants men black'
$ws.Cells.Item(99, 1).Value = 'football tights youth'
$ws.Cells.Item(100, 1).Value = 'mens 3/4 tights'
